# Rename existing sheet "Sheet1" -> "annual"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "annual"

# Add new sheet "qtr" after "annual"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "qtr"

# Fill in header row
$ws2.Range("A1").Value = "year"
$ws2.Range("B1").Value = "Q1"
$ws2.Range("C1").Value = "Q2"
$ws2.Range("D1").Value = "Q3"
$ws2.Range("E1").Value = "Q4"

# Fill in data rows (entered in the order the author typed them: the
# actual/forecast rows for 2023/2024 first, then the 2022f row inserted
# above them afterwards -- this matches the resulting sharedStrings order)
$ws2.Range("A3").Value = "2023a"
$ws2.Range("B3").Value = 9892
$ws2.Range("C3").Value = 9596
$ws2.Range("D3").Value = 9930
$ws2.Range("E3").Value = 10525

$ws2.Range("A4").Value = "2023f"
$ws2.Range("B4").Value = 11434
$ws2.Range("C4").Value = 11077
$ws2.Range("D4").Value = 11162
$ws2.Range("E4").Value = 11145

$ws2.Range("A5").Value = "2024a"
$ws2.Range("B5").Value = 9372
$ws2.Range("C5").Value = 9786
$ws2.Range("D5").Value = 10427

$ws2.Range("A2").Value = "2022f"
$ws2.Range("B2").Value = 9112
$ws2.Range("C2").Value = 8335
$ws2.Range("D2").Value = 10985
$ws2.Range("E2").Value = 11884

# Select A3 as active cell on qtr sheet
$ws2.Range("A3").Select()
